$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 60

$ws.Cells.Item($row, 1).Value = 5
$ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($row, 3).Value = "Maule"
$ws.Cells.Item($row, 4).Value = 44448
$ws.Cells.Item($row, 4).Style = $ws.Cells.Item(59, 4).Style
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(59, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 7
$ws.Cells.Item($row, 6).Value = 100112001
$ws.Cells.Item($row, 7).Value = "Berenjena"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 200
$ws.Cells.Item($row, 11).Value = 7000
$ws.Cells.Item($row, 12).Value = 7000
$ws.Cells.Item($row, 13).Value = 7000
$ws.Cells.Item($row, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 140
$ws.Cells.Item($row, 17).Value = 50
$ws.Cells.Item($row, 18).Value = "Hortaliza"
